# Edit script for NEBRASKA_2016.xlsx-style change:
# 1) Rename header labels (row 1) to short machine-friendly names
# 2) Capitalize the Spanish linking words (de/del/la/las/los/el) in
#    "Municipio Origen"/"Estado de Origen" place names
# 3) Tiny float re-serialization on D470 (rounding of a percentage)
# 4) Remove the trailing footnote rows (803-807) and shrink the used
#    range dimension back down to A1:D801

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Header renames ---------------------------------------------------
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- 2) Municipality / state name capitalization fixes -------------------
$ws.Range("B6").Value = "Pabellón De Arteaga"
$ws.Range("B7").Value = "Rincón De Romos"
$ws.Range("B23").Value = "Amatenango De La Frontera"
$ws.Range("B27").Value = "Comitán De Domínguez"
$ws.Range("B38").Value = "Ocozocoautla De Espinosa"
$ws.Range("B43").Value = "San Cristóbal De Las Casas"
$ws.Range("B71").Value = "Hidalgo Del Parral"
$ws.Range("B84").Value = "San Francisco De Borja"
$ws.Range("B85").Value = "San Francisco Del Oro"
$ws.Range("B88").Value = "Valle De Zaragoza"
$ws.Range("B109").Value = "Villa De Álvarez"
$ws.Range("A111").Value = "Ciudad De México"
$ws.Range("B126").Value = "Coneto De Comonfort"
$ws.Range("B140").Value = "Nombre De Dios"
$ws.Range("B142").Value = "Pánuco De Coronado"
$ws.Range("B148").Value = "San Juan De Guadalupe"
$ws.Range("B149").Value = "San Juan Del Río"
$ws.Range("B150").Value = "San Luis Del Cordero"
$ws.Range("A158").Value = "Estado De México"
$ws.Range("B158").Value = "Almoloya De Alquisiras"
$ws.Range("B161").Value = "Atizapán De Zaragoza"
$ws.Range("B165").Value = "Chapa De Mota"
$ws.Range("B167").Value = "Coacalco De Berriozábal"
$ws.Range("B170").Value = "Ecatepec De Morelos"
$ws.Range("B172").Value = "Ixtapan De La Sal"
$ws.Range("B177").Value = "Naucalpan De Juárez"
$ws.Range("B181").Value = "San Felipe Del Progreso"
$ws.Range("B189").Value = "Tlalnepantla De Baz"
$ws.Range("B193").Value = "Valle De Bravo"
$ws.Range("B194").Value = "Valle De Chalco Solidaridad"
$ws.Range("B195").Value = "Villa Del Carbón"
$ws.Range("B204").Value = "San Miguel De Allende"
$ws.Range("B205").Value = "Apaseo El Alto"
$ws.Range("B206").Value = "Apaseo El Grande"
$ws.Range("B212").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B216").Value = "Jaral Del Progreso"
$ws.Range("B223").Value = "Purísima Del Rincón"
$ws.Range("B228").Value = "San Francisco Del Rincón"
$ws.Range("B230").Value = "San Luis De La Paz"
$ws.Range("B231").Value = "Santa Cruz De Juventino Rosas"
$ws.Range("B232").Value = "Silao De La Victoria"
$ws.Range("B235").Value = "Valle De Santiago"
$ws.Range("B241").Value = "Acapulco De Juárez"
$ws.Range("B243").Value = "Ajuchitlán Del Progreso"
$ws.Range("B244").Value = "Alcozauca De Guerrero"
$ws.Range("B247").Value = "Atoyac De Álvarez"
$ws.Range("B248").Value = "Ayutla De Los Libres"
$ws.Range("B251").Value = "Buenavista De Cuéllar"
$ws.Range("B252").Value = "Chilapa De Álvarez"
$ws.Range("B253").Value = "Chilpancingo De Los Bravo"
$ws.Range("B254").Value = "Coahuayutla De José María Izazaga"
$ws.Range("B257").Value = "Coyuca De Benítez"
$ws.Range("B258").Value = "Coyuca De Catalán"
$ws.Range("B262").Value = "Cuetzala Del Progreso"
$ws.Range("B263").Value = "Cutzamala De Pinzón"
$ws.Range("B268").Value = "Huitzuco De Los Figueroa"
$ws.Range("B269").Value = "Iguala De La Independencia"
$ws.Range("B270").Value = "Zihuatanejo De Azueta"
$ws.Range("B282").Value = "Taxco De Alarcón"
$ws.Range("B284").Value = "Técpan De Galeana"
$ws.Range("B286").Value = "Tepecoacuilco De Trujano"
$ws.Range("B287").Value = "Tixtla De Guerrero"
$ws.Range("B289").Value = "Tlapa De Comonfort"
$ws.Range("B297").Value = "Atotonilco El Grande"
$ws.Range("B299").Value = "Cuautepec De Hinojosa"
$ws.Range("B301").Value = "Huasca De Ocampo"
$ws.Range("B304").Value = "Jacala De Ledezma"
$ws.Range("B306").Value = "Mixquiahuala De Juárez"
$ws.Range("B307").Value = "Nopala De Villagrán"
$ws.Range("B308").Value = "Pachuca De Soto"
$ws.Range("B310").Value = "Progreso De Obregón"
$ws.Range("B312").Value = "Tepeji Del Río De Ocampo"
$ws.Range("B316").Value = "Tulancingo De Bravo"
$ws.Range("B325").Value = "Atotonilco El Alto"
$ws.Range("B326").Value = "Autlán De Navarro"
$ws.Range("B340").Value = "Encarnación De Díaz"
$ws.Range("B349").Value = "Jilotlán De Los Dolores"
$ws.Range("B354").Value = "La Manzanilla De La Paz"
$ws.Range("B355").Value = "Lagos De Moreno"
$ws.Range("B360").Value = "Ojuelos De Jalisco"
$ws.Range("B365").Value = "San Diego De Alejandría"
$ws.Range("B366").Value = "San Juan De Los Lagos"
$ws.Range("B369").Value = "San Miguel El Alto"
$ws.Range("B371").Value = "Tamazula De Gordiano"
$ws.Range("B375").Value = "Tepatitlán De Morelos"
$ws.Range("B376").Value = "Tizapán El Alto"
$ws.Range("B377").Value = "Tlajomulco De Zúñiga"
$ws.Range("B385").Value = "Unión De San Antonio"
$ws.Range("B386").Value = "Unión De Tula"
$ws.Range("B390").Value = "Yahualica De González Gallo"
$ws.Range("B393").Value = "Zapotlán Del Rey"
$ws.Range("B480").Value = "Puente De Ixtla"
$ws.Range("B483").Value = "Tlaltizapán De Zapata"
$ws.Range("B494").Value = "Ixtlán Del Río"
$ws.Range("B499").Value = "Santa María Del Oro"
$ws.Range("B515").Value = "Acatlán De Pérez Figueroa"
$ws.Range("B518").Value = "Ayoquezco De Aldama"
$ws.Range("B520").Value = "Chalcatongo De Hidalgo"
$ws.Range("B522").Value = "Constancia Del Rosario"
$ws.Range("B524").Value = "Heroica Ciudad De Ejutla De Crespo"
$ws.Range("B525").Value = "Heroica Ciudad De Huajuapan De León"
$ws.Range("B526").Value = "Heroica Ciudad De Tlaxiaco"
$ws.Range("B528").Value = "Ixtlán De Juárez"
$ws.Range("B529").Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Range("B531").Value = "Mariscala De Juárez"
$ws.Range("B534").Value = "Miahuatlán De Porfirio Díaz"
$ws.Range("B535").Value = "Oaxaca De Juárez"
$ws.Range("B536").Value = "Putla Villa De Guerrero"
$ws.Range("B546").Value = "San Juan Del Estado"
$ws.Range("B561").Value = "Santa Cruz Tacache De Mina"
$ws.Range("B576").Value = "Tlacolula De Matamoros"
$ws.Range("B577").Value = "Villa De Tututepec"
$ws.Range("B578").Value = "Villa Sola De Vega"
$ws.Range("B580").Value = "Zimatlán De Álvarez"
$ws.Range("B596").Value = "Huehuetlán El Chico"
$ws.Range("B598").Value = "Izúcar De Matamoros"
$ws.Range("B601").Value = "Los Reyes De Juárez"
$ws.Range("B613").Value = "Tlacotepec De Benito Juárez"
$ws.Range("B621").Value = "Amealco De Bonfil"
$ws.Range("B623").Value = "Cadereyta De Montes"
$ws.Range("B627").Value = "Jalpan De Serra"
$ws.Range("B628").Value = "Landa De Matamoros"
$ws.Range("B632").Value = "San Juan Del Río"
$ws.Range("B641").Value = "Ciudad Del Maíz"
$ws.Range("B646").Value = "Mexquitic De Carmona"
$ws.Range("B651").Value = "Santa María Del Río"
$ws.Range("B654").Value = "Villa De Arriaga"
$ws.Range("B655").Value = "Villa De Ramos"
$ws.Range("B656").Value = "Villa De Reyes"
$ws.Range("B713").Value = "Alto Lucero De Gutiérrez Barrios"
$ws.Range("B717").Value = "Castillo De Teayo"
$ws.Range("B727").Value = "Hueyapan De Ocampo"
$ws.Range("B728").Value = "Ignacio De La Llave"
$ws.Range("B730").Value = "Ixhuatlán Del Café"
$ws.Range("B735").Value = "Juchique De Ferrer"
$ws.Range("B738").Value = "Martínez De La Torre"
$ws.Range("B745").Value = "Paso De Ovejas"
$ws.Range("B749").Value = "Soledad De Doblado"
$ws.Range("B792").Value = "Tlaltenango De Sánchez Román"
$ws.Range("B795").Value = "Villa De Cos"

# --- 3) Tiny floating point re-serialization on D470 ----------------------
$ws.Range("D470").Value = 0.0933237616654702

# --- 4) Drop the trailing footnote rows (803-807) -------------------------
$ws.Range("A803:A807").EntireRow.Delete()
